$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current row 184, pushing the existing
# rows 184:207 down to 189:212.
$ws.Rows("184:188").Insert()

# Common values shared by every row in this block.
$mercado = "Agrícola del Norte S.A. de Arica"
$region  = "Arica y Parinacota"
$codreg  = 15
$tipo    = "Fruta"
$prodId  = 100104
$prod    = "Frutos de pepita"
$catId   = 100104002
$cat     = "Manzana"
$origen  = "Región de O'Higgins"
$unidad  = "`$/caja 16 kilos empedrada"
$fecha   = 45034

# New rows' variety-specific data:
# row, variedad, calidad, volumen, precioMin, precioMax, precioProm, precioKg, kgUnidad
$newRows = @(
  @(184, "Ambrosia",            "Segunda", 320, 20000, 22000, 21125, 1320, 16),
  @(185, "Fuji royal",          "Segunda", 300, 20000, 22000, 21000, 1312, 16),
  @(186, "Granny Smith",        "Segunda", 300, 20000, 22000, 21000, 1312, 16),
  @(187, "Richared Delicious",  "Segunda", 300, 20000, 22000, 21000, 1312, 16),
  @(188, "Royal Gala",          "Segunda", 350, 20000, 22000, 21143, 1321, 16)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $prodId
    $ws.Cells.Item($r, 8).Value = $prod
    $ws.Cells.Item($r, 9).Value = $catId
    $ws.Cells.Item($r, 10).Value = $cat
    $ws.Cells.Item($r, 11).Value = $row[1]
    $ws.Cells.Item($r, 12).Value = $row[2]
    $ws.Cells.Item($r, 13).Value = $row[3]
    $ws.Cells.Item($r, 14).Value = $row[4]
    $ws.Cells.Item($r, 15).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $row[7]
    $ws.Cells.Item($r, 20).Value = $row[8]
}

$ws.Range("A1").Select() | Out-Null
